$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Table layout: col A = index, B = name, C = from_bus, D = to_bus, E = in_service
# Two new contingency lines (line7, line8) are inserted into the result set
# right after line6, which bumps every "extr*" case down by two rows and
# refreshes their computed bus values; two brand-new rows are appended at
# the bottom for the final two extra cases.

# Row 8: was extr1 -> becomes line7
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11

# Row 9: was extr2 -> becomes line8
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 5).Value = $true

# Row 10: was extr3 -> becomes extr1
$ws.Cells.Item(10, 2).Value = "extr1"
$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = 12
$ws.Cells.Item(10, 5).Value = $true

# Row 11: was extr4 -> becomes extr2
$ws.Cells.Item(11, 2).Value = "extr2"
$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 9
$ws.Cells.Item(11, 5).Value = $true

# Row 12: was extr5 -> becomes extr3
$ws.Cells.Item(12, 2).Value = "extr3"
$ws.Cells.Item(12, 3).Value = 10

# Row 13: was extr6 -> becomes extr4
$ws.Cells.Item(13, 2).Value = "extr4"
$ws.Cells.Item(13, 4).Value = 8

# Row 14: was extr7 -> becomes extr5
$ws.Cells.Item(14, 2).Value = "extr5"
$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 11

# Row 15: was extr8 -> becomes extr6
$ws.Cells.Item(15, 2).Value = "extr6"
$ws.Cells.Item(15, 3).Value = 7
$ws.Cells.Item(15, 4).Value = 11
$ws.Cells.Item(15, 5).Value = $false

# New row 16: extr7
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 1).Font.Bold = $true
$ws.Cells.Item(16, 1).HorizontalAlignment = -4108
$ws.Cells.Item(16, 1).VerticalAlignment = -4160
$ws.Cells.Item(16, 1).Borders.LineStyle = 1
$ws.Cells.Item(16, 2).Value = "extr7"
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = $false

# New row 17: extr8
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 1).Font.Bold = $true
$ws.Cells.Item(17, 1).HorizontalAlignment = -4108
$ws.Cells.Item(17, 1).VerticalAlignment = -4160
$ws.Cells.Item(17, 1).Borders.LineStyle = 1
$ws.Cells.Item(17, 2).Value = "extr8"
$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = $true
